# Applies "RMI updates through 12/2" edits to BAU Subsidies.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the Calculations sheet: "Petroleum Products" source note + data
# ---------------------------------------------------------------------------
$calc = $wb.Worksheets.Item("Calculations")

# Row 4 ("Capacity (MW)" style block) - update source note and yearly values
$calc.Range("B4").Value = "See elec/CCaMC"

$row4Values = @(
    1462790.2009077901, 1284828.0428722501, 1237544.2340474708, 1190260.4252226916,
    1142976.6163979124, 1095692.8075731331, 1048408.9987483537, 1001125.1899235747,
    953841.38109879522, 906557.57227401587, 859273.76344923663, 811989.95462445728,
    764706.14579967817, 757913.21727074252, 751120.28874180699, 744327.36021287122,
    737534.43168393557, 730741.50315499993, 723948.57462606428, 717155.64609712863,
    710362.71756819298, 703569.78903925733, 696776.86051032168, 689983.93198138592,
    683191.00345245027, 676398.07492351462, 669605.14639457897, 662812.21786564332,
    656019.28933670768, 649226.36080777203, 642433.43227883638, 635640.50374990073
)
for ($i = 0; $i -lt $row4Values.Length; $i++) {
    $calc.Cells.Item(4, 3 + $i).Value = $row4Values[$i]
}

# Row 10 - update source note and yearly values
$calc.Range("B10").Value = "See elec/CCaMC"

$row10Values = @(
    6831836.4795198459, 6500515.8491650894, 6169195.2188103329, 5843882.7444435377,
    5630238.7966989167, 5422794.6290817559, 5233450.6532794116, 5058216.8964543967,
    4898609.074266511, 4753976.9239266422, 4620219.8524136795, 4501029.6773458235,
    4393594.5403691512, 4297491.1596871642, 4212533.3580054678, 4135494.7609418505,
    4069499.8884631144, 4010225.5494396384, 3958592.5068339193, 3914721.5043390612,
    3876446.5635061474, 3843367.6747004823, 3814239.6224356443, 3789350.4384214408,
    3768365.5466244677, 3749025.879948609, 3730662.2863270584, 3714585.4148475262,
    3697197.8942802823, 3680144.9497464434, 3662115.9321582974, 3640824.119373824
)
for ($i = 0; $i -lt $row10Values.Length; $i++) {
    $calc.Cells.Item(10, 3 + $i).Value = $row10Values[$i]
}

$wb.Save()
